$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 23.76450733333333
$ws.Range("H2").Value = 71.293522
$ws.Range("I2").Value = 0.3430552013751479
$ws.Range("J2").Value = 0.3430552013751479
$ws.Range("M2").Value = 6.111751666666666
$ws.Range("N2").Value = 18.335255
$ws.Range("O2").Value = 0.6061514841909396
$ws.Range("P2").Value = 0.6061514841909394
$ws.Range("Q2").Value = 145.2427673020122
$ws.Range("R2").Value = 1307.18490571811
$ws.Range("S2").Value = 0.2079434194729675
$ws.Range("T2").Value = 0.2079434194729675
$ws.Range("G3").Value = 23.76450733333333
$ws.Range("H3").Value = 71.293522
$ws.Range("I3").Value = 0.3430552013751479
$ws.Range("J3").Value = 0.3430552013751479
$ws.Range("O3").Value = 0.2731664420559804
$ws.Range("P3").Value = 0.2731664420559804
$ws.Range("Q3").Value = 65.45467760622932
$ws.Range("R3").Value = 589.0920984560639
$ws.Range("S3").Value = 0.09371116878844703
$ws.Range("T3").Value = 0.09371116878844703
$ws.Range("G4").Value = 23.76450733333333
$ws.Range("H4").Value = 71.293522
$ws.Range("I4").Value = 0.3430552013751479
$ws.Range("J4").Value = 0.3430552013751479
$ws.Range("M4").Value = 0.568439
$ws.Range("N4").Value = 1.705317
$ws.Range("O4").Value = 0.0563766596410053
$ws.Range("P4").Value = 0.05637665964100529
$ws.Range("Q4").Value = 13.50867278405267
$ws.Range("R4").Value = 121.578055056474
$ws.Range("S4").Value = 0.01934030632600325
$ws.Range("T4").Value = 0.01934030632600324
$ws.Range("G5").Value = 23.76450733333333
$ws.Range("H5").Value = 71.293522
$ws.Range("I5").Value = 0.3430552013751479
$ws.Range("J5").Value = 0.3430552013751479
$ws.Range("M5").Value = 0.3689163333333333
$ws.Range("N5").Value = 1.106749
$ws.Range("O5").Value = 0.036588394815171
$ws.Range("P5").Value = 0.036588394815171
$ws.Range("Q5").Value = 8.767114908886445
$ws.Range("R5").Value = 78.904034179978
$ws.Range("S5").Value = 0.01255183915131191
$ws.Range("T5").Value = 0.01255183915131191
$ws.Range("G6").Value = 23.76450733333333
$ws.Range("H6").Value = 71.293522
$ws.Range("I6").Value = 0.3430552013751479
$ws.Range("J6").Value = 0.3430552013751479
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2794673333333333
$ws.Range("N6").Value = 0.838402
$ws.Range("O6").Value = 0.02771701929690382
$ws.Range("P6").Value = 0.02771701929690381
$ws.Range("Q6").Value = 6.641403492427111
$ws.Range("R6").Value = 59.772631431844
$ws.Range("S6").Value = 0.009508467636418199
$ws.Range("T6").Value = 0.009508467636418197
$ws.Range("I7").Value = 0.3627390578207265
$ws.Range("J7").Value = 0.3627390578207265
$ws.Range("M7").Value = 6.111751666666666
$ws.Range("N7").Value = 18.335255
$ws.Range("O7").Value = 0.6061514841909396
$ws.Range("P7").Value = 0.6061514841909394
$ws.Range("Q7").Value = 153.5765216653661
$ws.Range("R7").Value = 1382.188694988295
$ws.Range("S7").Value = 0.2198748182720564
$ws.Range("T7").Value = 0.2198748182720564
$ws.Range("I8").Value = 0.3627390578207265
$ws.Range("J8").Value = 0.3627390578207265
$ws.Range("O8").Value = 0.2731664420559804
$ws.Range("P8").Value = 0.2731664420559804
$ws.Range("S8").Value = 0.09908813781962642
$ws.Range("T8").Value = 0.09908813781962642
$ws.Range("I9").Value = 0.3627390578207265
$ws.Range("J9").Value = 0.3627390578207265
$ws.Range("M9").Value = 0.568439
$ws.Range("N9").Value = 1.705317
$ws.Range("O9").Value = 0.0563766596410053
$ws.Range("P9").Value = 0.05637665964100529
$ws.Range("Q9").Value = 14.28377479325033
$ws.Range("R9").Value = 128.553973139253
$ws.Range("S9").Value = 0.02045001640125804
$ws.Range("T9").Value = 0.02045001640125804
$ws.Range("I10").Value = 0.3627390578207265
$ws.Range("J10").Value = 0.3627390578207265
$ws.Range("M10").Value = 0.3689163333333333
$ws.Range("N10").Value = 1.106749
$ws.Range("O10").Value = 0.036588394815171
$ws.Range("P10").Value = 0.036588394815171
$ws.Range("Q10").Value = 9.270155325171222
$ws.Range("R10").Value = 83.431397926541
$ws.Range("S10").Value = 0.01327203986242789
$ws.Range("T10").Value = 0.01327203986242789
$ws.Range("I11").Value = 0.3627390578207265
$ws.Range("J11").Value = 0.3627390578207265
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2794673333333333
$ws.Range("N11").Value = 0.838402
$ws.Range("O11").Value = 0.02771701929690382
$ws.Range("P11").Value = 0.02771701929690381
$ws.Range("Q11").Value = 7.022474621557556
$ws.Range("R11").Value = 63.202271594018
$ws.Range("S11").Value = 0.01005404546535779
$ws.Range("T11").Value = 0.01005404546535778
$ws.Range("G12").Value = 2.481784666666667
$ws.Range("H12").Value = 7.445354
$ws.Range("I12").Value = 0.03582607990357473
$ws.Range("J12").Value = 0.03582607990357473
$ws.Range("M12").Value = 6.111751666666666
$ws.Range("N12").Value = 18.335255
$ws.Range("O12").Value = 0.6061514841909396
$ws.Range("P12").Value = 0.6061514841909394
$ws.Range("Q12").Value = 15.16805157280777
$ws.Range("R12").Value = 136.51246415527
$ws.Range("S12").Value = 0.02171603150629502
$ws.Range("T12").Value = 0.02171603150629501
$ws.Range("G13").Value = 2.481784666666667
$ws.Range("H13").Value = 7.445354
$ws.Range("I13").Value = 0.03582607990357473
$ws.Range("J13").Value = 0.03582607990357473
$ws.Range("O13").Value = 0.2731664420559804
$ws.Range("P13").Value = 0.2731664420559804
$ws.Range("Q13").Value = 6.835589434538666
$ws.Range("R13").Value = 61.520304910848
$ws.Range("S13").Value = 0.009786482780072771
$ws.Range("T13").Value = 0.009786482780072771
$ws.Range("G14").Value = 2.481784666666667
$ws.Range("H14").Value = 7.445354
$ws.Range("I14").Value = 0.03582607990357473
$ws.Range("J14").Value = 0.03582607990357473
$ws.Range("M14").Value = 0.568439
$ws.Range("N14").Value = 1.705317
$ws.Range("O14").Value = 0.0563766596410053
$ws.Range("P14").Value = 0.05637665964100529
$ws.Range("Q14").Value = 1.410743194135333
$ws.Range("R14").Value = 12.696688747218
$ws.Range("S14").Value = 0.002019754712995292
$ws.Range("T14").Value = 0.002019754712995292
$ws.Range("G15").Value = 2.481784666666667
$ws.Range("H15").Value = 7.445354
$ws.Range("I15").Value = 0.03582607990357473
$ws.Range("J15").Value = 0.03582607990357473
$ws.Range("M15").Value = 0.3689163333333333
$ws.Range("N15").Value = 1.106749
$ws.Range("O15").Value = 0.036588394815171
$ws.Range("P15").Value = 0.036588394815171
$ws.Range("Q15").Value = 0.9155708993495555
$ws.Range("R15").Value = 8.240138094145999
$ws.Range("S15").Value = 0.001310818756191856
$ws.Range("T15").Value = 0.001310818756191856
$ws.Range("G16").Value = 2.481784666666667
$ws.Range("H16").Value = 7.445354
$ws.Range("I16").Value = 0.03582607990357473
$ws.Range("J16").Value = 0.03582607990357473
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2794673333333333
$ws.Range("N16").Value = 0.838402
$ws.Range("O16").Value = 0.02771701929690382
$ws.Range("P16").Value = 0.02771701929690381
$ws.Range("Q16").Value = 0.6935777427008889
$ws.Range("R16").Value = 6.242199684308
$ws.Range("S16").Value = 0.0009929921480197989
$ws.Range("T16").Value = 0.0009929921480197987
$ws.Range("G17").Value = 9.395678666666667
$ws.Range("H17").Value = 28.187036
$ws.Range("I17").Value = 0.1356323693918298
$ws.Range("J17").Value = 0.1356323693918298
$ws.Range("M17").Value = 6.111751666666666
$ws.Range("N17").Value = 18.335255
$ws.Range("O17").Value = 0.6061514841909396
$ws.Range("P17").Value = 0.6061514841909394
$ws.Range("Q17").Value = 57.42405475046444
$ws.Range("R17").Value = 516.81649275418
$ws.Range("S17").Value = 0.08221376201119139
$ws.Range("T17").Value = 0.08221376201119138
$ws.Range("G18").Value = 9.395678666666667
$ws.Range("H18").Value = 28.187036
$ws.Range("I18").Value = 0.1356323693918298
$ws.Range("J18").Value = 0.1356323693918298
$ws.Range("O18").Value = 0.2731664420559804
$ws.Range("P18").Value = 0.2731664420559804
$ws.Range("Q18").Value = 25.87855533431467
$ws.Range("R18").Value = 232.906998008832
$ws.Range("S18").Value = 0.0370502117743886
$ws.Range("T18").Value = 0.0370502117743886
$ws.Range("G19").Value = 9.395678666666667
$ws.Range("H19").Value = 28.187036
$ws.Range("I19").Value = 0.1356323693918298
$ws.Range("J19").Value = 0.1356323693918298
$ws.Range("M19").Value = 0.568439
$ws.Range("N19").Value = 1.705317
$ws.Range("O19").Value = 0.0563766596410053
$ws.Range("P19").Value = 0.05637665964100529
$ws.Range("Q19").Value = 5.340870185601334
$ws.Range("R19").Value = 48.067831670412
$ws.Range("S19").Value = 0.007646499925506293
$ws.Range("T19").Value = 0.007646499925506292
$ws.Range("G20").Value = 9.395678666666667
$ws.Range("H20").Value = 28.187036
$ws.Range("I20").Value = 0.1356323693918298
$ws.Range("J20").Value = 0.1356323693918298
$ws.Range("M20").Value = 0.3689163333333333
$ws.Range("N20").Value = 1.106749
$ws.Range("O20").Value = 0.036588394815171
$ws.Range("P20").Value = 0.036588394815171
$ws.Range("Q20").Value = 3.466219322884889
$ws.Range("R20").Value = 31.195973905964
$ws.Range("S20").Value = 0.004962570681025383
$ws.Range("T20").Value = 0.004962570681025383
$ws.Range("G21").Value = 9.395678666666667
$ws.Range("H21").Value = 28.187036
$ws.Range("I21").Value = 0.1356323693918298
$ws.Range("J21").Value = 0.1356323693918298
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2794673333333333
$ws.Range("N21").Value = 0.838402
$ws.Range("O21").Value = 0.02771701929690382
$ws.Range("P21").Value = 0.02771701929690381
$ws.Range("Q21").Value = 2.625785261830222
$ws.Range("R21").Value = 23.632067356472
$ws.Range("S21").Value = 0.003759324999718133
$ws.Range("T21").Value = 0.003759324999718133
$ws.Range("G22").Value = 8.503089000000001
$ws.Range("H22").Value = 25.509267
$ws.Range("I22").Value = 0.122747291508721
$ws.Range("J22").Value = 0.122747291508721
$ws.Range("M22").Value = 6.111751666666666
$ws.Range("N22").Value = 18.335255
$ws.Range("O22").Value = 0.6061514841909396
$ws.Range("P22").Value = 0.6061514841909394
$ws.Range("Q22").Value = 51.968768367565
$ws.Range("R22").Value = 467.7189153080849
$ws.Range("S22").Value = 0.07440345292842918
$ws.Range("T22").Value = 0.07440345292842915
$ws.Range("G23").Value = 8.503089000000001
$ws.Range("H23").Value = 25.509267
$ws.Range("I23").Value = 0.122747291508721
$ws.Range("J23").Value = 0.122747291508721
$ws.Range("O23").Value = 0.2731664420559804
$ws.Range("P23").Value = 0.2731664420559804
$ws.Range("Q23").Value = 23.420092045056
$ws.Range("R23").Value = 210.780828405504
$ws.Range("S23").Value = 0.03353044089344558
$ws.Range("T23").Value = 0.03353044089344558
$ws.Range("G24").Value = 8.503089000000001
$ws.Range("H24").Value = 25.509267
$ws.Range("I24").Value = 0.122747291508721
$ws.Range("J24").Value = 0.122747291508721
$ws.Range("M24").Value = 0.568439
$ws.Range("N24").Value = 1.705317
$ws.Range("O24").Value = 0.0563766596410053
$ws.Range("P24").Value = 0.05637665964100529
$ws.Range("Q24").Value = 4.833487408071001
$ws.Range("R24").Value = 43.501386672639
$ws.Range("S24").Value = 0.006920082275242426
$ws.Range("T24").Value = 0.006920082275242424
$ws.Range("G25").Value = 8.503089000000001
$ws.Range("H25").Value = 25.509267
$ws.Range("I25").Value = 0.122747291508721
$ws.Range("J25").Value = 0.122747291508721
$ws.Range("M25").Value = 0.3689163333333333
$ws.Range("N25").Value = 1.106749
$ws.Range("O25").Value = 0.036588394815171
$ws.Range("P25").Value = 0.036588394815171
$ws.Range("Q25").Value = 3.136928415887001
$ws.Range("R25").Value = 28.232355742983
$ws.Range("S25").Value = 0.004491126364213973
$ws.Range("T25").Value = 0.004491126364213972
$ws.Range("G26").Value = 8.503089000000001
$ws.Range("H26").Value = 25.509267
$ws.Range("I26").Value = 0.122747291508721
$ws.Range("J26").Value = 0.122747291508721
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2794673333333333
$ws.Range("N26").Value = 0.838402
$ws.Range("O26").Value = 0.02771701929690382
$ws.Range("P26").Value = 0.02771701929690381
$ws.Range("Q26").Value = 2.376335607926
$ws.Range("R26").Value = 21.387020471334
$ws.Range("S26").Value = 0.003402189047389899
$ws.Range("T26").Value = 0.003402189047389898
